$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A26:C26").Copy()
$ws.Range("A27:C27").PasteSpecial(-4122)

$ws.Range("A27").Value = 45267
$ws.Range("B27").Value = "Internship"
$ws.Range("C27").Value = "Contributed technical work by aiding in resolving inconsistencies flagged by the system for employee calls"

$ws.Range("C28").Select()
